# Applies:
#  1. Delete slide 1 ("Colour scheme") - only the "DCM model images" slide remains.
#  2. Update the cached datetimeFigureOut placeholder text ("28/11/2016" ->
#     "04/05/2017") on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- Remove the "Colour scheme" slide (keeps only the DCM model images slide) ---
$p.Slides.Item(1).Delete()

$oldDate = "28/11/2016"
$newDate = "04/05/2017"
$ppPlaceholderDate = 16

# --- Refresh the cached "last printed" date field on the slide master ---
$master = $p.SlideMaster
$masterShapes = $master.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Refresh the cached "last printed" date field on every slide layout ---
for ($k = 1; $k -le $master.CustomLayouts.Count; $k++) {
    $layout = $master.CustomLayouts.Item($k)
    $layoutShapes = $layout.Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $sh2 = $layoutShapes.Item($j)
        if ($sh2.Type -eq 14 -and $sh2.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($sh2.TextFrame.TextRange.Text -eq $oldDate) {
                $sh2.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
